# Weekly update: insert a new daily price record for "Rabanito" at
# Vega Modelo de Temuco, pushing the existing historical rows (32-118)
# down by one row (to 33-119).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 32, shifting rows
# 32..118 down to 33..119 (Excel also extends the dimension to R119
# and carries the date-column style (s="2") down automatically).
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record.
$ws.Range("A32").Value = 10
$ws.Range("B32").Value = "Vega Modelo de Temuco"
$ws.Range("C32").Value = "La Araucanía"
$ws.Range("D32").Value = 45114
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = 300000001
$ws.Range("G32").Value = "Rabanito"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 70
$ws.Range("K32").Value = 7000
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = 7429
$ws.Range("N32").Value = "$/docena de paquetes"
$ws.Range("O32").Value = "Provincia de Cautín"
$ws.Range("P32").Value = 619
$ws.Range("Q32").Value = 12
$ws.Range("R32").Value = "Hortaliza"
